$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the first "# ACTION POINTS Week4 #" title paragraph
#    (the one that also carries the OLE_LINK1 bookmark).
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$clearRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$clearRange.Text = ""
$d.Paragraphs.Item(2).Range.Delete()

# ------------------------------------------------------------------
# 2. Turn the "Fix PIP install ( Marco ) shame on you continue"
#    numbered-list paragraph into the new Title paragraph
#    "# ACTION POINTS #".
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Style = "Title"
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = "# ACTION POINTS #"

# ------------------------------------------------------------------
# 3. Drop the now-redundant second title block: the trailing blank
#    paragraph, the second "# ACTION POINTS Week4 #" title, and the
#    blank paragraph that used to sit right after the first title.
#    Deleted highest index first so earlier indices stay valid.
# ------------------------------------------------------------------
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(6).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()

# ------------------------------------------------------------------
# 4. Text-only edits on the remaining action-item paragraphs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Complete Exercise 3 file (Iryna)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Complete Exercise 4 file (Iryna)", 2) | Out-Null

$d.Content.Find.Execute("Upload calculator scrip ( Iryna )", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Study Lessons Files:", 2) | Out-Null

$d.Content.Find.Execute("Study Modules and Loops", $true, $false, $false, $false, $false, `
    $true, 1, $false, "textFiles", 2) | Out-Null

$d.Content.Find.Execute("Complete debugging Exercise (Iryna)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "JsonFiles", 2) | Out-Null

# ------------------------------------------------------------------
# 5. Add a new "CVVFiles" line using the trailing blank paragraph.
# ------------------------------------------------------------------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastP.Range.Text = "CVVFiles"
$d.Paragraphs.Item($d.Paragraphs.Count).Range.LanguageID = "en-US"
